# "separated DT from RFC"
# The Decision Tree model rating (B8) was nudged from 82.3 to 82.2 so that it
# no longer shares the same rounded score as the Random Forest Classifier.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Decision Tree rating value.
$ws.Range("B8").Value = 82.2

# Move the active selection to B11, matching where the cursor ended up
# after the edit.
$ws.Range("B11").Select()
